$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (column J) mirroring the formatting of column I
# for rows 4 (header) through 14 (totals row), with the new data values.

$values = @{
    4  = 2021
    5  = 24.4
    6  = 45.7
    7  = 38
    8  = 51.3
    9  = 51.5
    10 = 13
    11 = 36.4
    12 = 27
    13 = 2.7
    14 = 40.4
}

foreach ($row in 4..14) {
    $srcCell = $ws.Cells.Item($row, 9)   # column I
    $dstCell = $ws.Cells.Item($row, 10)  # column J

    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)         # xlPasteFormats
    $dstCell.Value = $values[$row]
}

# Row 3 separator height changes from 18 to 13.5
$ws.Rows.Item(3).RowHeight = 13.5

# Update the selected cell shown in the saved view
$ws.Range("K18").Select()
